# Fix typo in Code Type values: "c2saferrust" -> "c2saferust"
# and "c2saferrustv2" -> "c2saferustv2" (one "r" removed from "saferr")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 9-15 currently hold "c2saferrust" in column A
for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "c2saferust"
}

# Rows 16-22 currently hold "c2saferrustv2" in column A
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = "c2saferustv2"
}
